$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "F:(1, 132)249.9, 0"
$ws.Range("C2").Value = "F:(1, 132)275.7, 0"
$ws.Range("D2").Value = "F:(1, 193)322.8, 0"

$ws.Range("B3").Value = "F:(2, 132)105.8, 0"
$ws.Range("C3").Value = "F:(2, 132)43.9, 0"
$ws.Range("D3").Value = "F:(3, 193)157.2, 0"

$ws.Range("B4").Value = "F:(4, 47)10.2, 0"
$ws.Range("C4").Value = "F:(4, 47)1.1, 0.3644"
$ws.Range("D4").Value = "F:(4, 46)2.4, 0.0605"

$ws.Range("B5").Value = "F:(8, 132)10.1, 0"
$ws.Range("C5").Value = "F:(8, 132)2.4, 0.0174"
$ws.Range("D5").Value = "F:(12, 193)2.2, 0.0123"
